$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1921.4286
$ws.Range("I112").Value = 1250
$ws.Range("J112").Value = 2190
$ws.Range("K112").Value = 3750
$ws.Range("L112").Value = 6570
$ws.Range("M112").Value = -2642
$ws.Range("N112").Value = -8786

$ws.Range("H118").Value = 4182.7856
$ws.Range("I118").Value = 4182
$ws.Range("K118").Value = 12546
$ws.Range("M118").Value = -10889

$ws.Range("H121").Value = 645.37036
$ws.Range("J121").Value = 633
$ws.Range("L121").Value = 1899
$ws.Range("N121").Value = -5393

$ws.Range("H132").Value = 1531.9807
$ws.Range("I132").Value = 1182.4524
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3547.357199999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1017.357199999999
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 50396.668
$ws.Range("I2").Value = 60719.41
$ws.Range("J2").Value = 6525
$ws.Range("K2").Value = 60719.41
$ws.Range("L2").Value = 6525
$ws.Range("M2").Value = -60606.41
$ws.Range("N2").Value = -6751

$ws.Range("H74").Value = 818.8889
$ws.Range("I74").Value = 811.0625
$ws.Range("K74").Value = 811.0625
$ws.Range("M74").Value = 62.9375

$ws.Range("H77").Value = 818.8889
$ws.Range("I77").Value = 811.0625
$ws.Range("K77").Value = 4055.3125
$ws.Range("M77").Value = 312.6875

$ws.Range("H116").Value = 50396.668
$ws.Range("I116").Value = 60719.41
$ws.Range("J116").Value = 6525
$ws.Range("K116").Value = 60719.41
$ws.Range("L116").Value = 6525
$ws.Range("M116").Value = -58425.41
$ws.Range("N116").Value = -11113

$ws.Range("H122").Value = 1021.88
$ws.Range("I122").Value = 960.2917
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 2880.8751
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -430.8751000000002
$ws.Range("N122").Value = -12400

$ws.Range("H132").Value = 1483.1052
$ws.Range("I132").Value = 1003.18604
$ws.Range("J132").Value = 2957.1428
$ws.Range("K132").Value = 3009.55812
$ws.Range("L132").Value = 8871.428400000001
$ws.Range("M132").Value = -479.5581200000001
$ws.Range("N132").Value = -13931.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 50396.668
$ws.Range("I3").Value = 60719.41
$ws.Range("J3").Value = 6525
$ws.Range("K3").Value = 60719.41
$ws.Range("L3").Value = 6525
$ws.Range("M3").Value = -60605.41
$ws.Range("N3").Value = -6753

$ws.Range("H26").Value = 6520.5
$ws.Range("I26").Value = 6520.5
$ws.Range("K26").Value = 6520.5
$ws.Range("M26").Value = -6228.5

$ws.Range("H86").Value = 2365.6155
$ws.Range("J86").Value = 2540.3
$ws.Range("L86").Value = 2540.3
$ws.Range("N86").Value = -4786.3

$ws.Range("H89").Value = 2365.6155
$ws.Range("J89").Value = 2540.3
$ws.Range("L89").Value = 12701.5
$ws.Range("N89").Value = -23933.5

$ws.Range("H94").Value = 1303.6
$ws.Range("I94").Value = 1233.25
$ws.Range("J94").Value = 1585
$ws.Range("K94").Value = 1233.25
$ws.Range("L94").Value = 1585
$ws.Range("M94").Value = -782.25
$ws.Range("N94").Value = -2487

$ws.Range("H96").Value = 9290.25
$ws.Range("I96").Value = 2563.3333
$ws.Range("K96").Value = 2563.3333
$ws.Range("M96").Value = 182.6667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1567.186
$ws.Range("I31").Value = 879.125
$ws.Range("J31").Value = 2436.3157
$ws.Range("K31").Value = 879.125
$ws.Range("L31").Value = 2436.3157
$ws.Range("M31").Value = -584.125
$ws.Range("N31").Value = -3026.3157

$ws.Range("H34").Value = 1567.186
$ws.Range("I34").Value = 879.125
$ws.Range("J34").Value = 2436.3157
$ws.Range("K34").Value = 879.125
$ws.Range("L34").Value = 2436.3157
$ws.Range("M34").Value = -677.125
$ws.Range("N34").Value = -2840.3157

$ws.Range("I122").Value = 909835.5600000001
$ws.Range("J122").Value = 710.5
$ws.Range("K122").Value = 2729506.68
$ws.Range("L122").Value = 2131.5
$ws.Range("M122").Value = -2727056.68
$ws.Range("N122").Value = -7031.5

$ws.Range("H132").Value = 2744.8696
$ws.Range("I132").Value = 1316.7142
$ws.Range("K132").Value = 3950.1426
$ws.Range("M132").Value = -1420.1426

$ws.Range("H140").Value = 52247.5
$ws.Range("J140").Value = 52247.5
$ws.Range("L140").Value = 52247.5
$ws.Range("N140").Value = -62607.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 117
$ws.Range("I40").Value = 122.625
$ws.Range("J40").Value = 94.5
$ws.Range("K40").Value = 490.5
$ws.Range("L40").Value = 378
$ws.Range("M40").Value = -421.5
$ws.Range("N40").Value = -516

$ws.Range("H92").Value = 50392.832
$ws.Range("I92").Value = 150095.5
$ws.Range("J92").Value = 541.5
$ws.Range("K92").Value = 450286.5
$ws.Range("L92").Value = 1624.5
$ws.Range("M92").Value = -449038.5
$ws.Range("N92").Value = -4120.5

$ws.Range("H96").Value = 3941.647
$ws.Range("J96").Value = 4063
$ws.Range("L96").Value = 12189
$ws.Range("N96").Value = -16307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1541.5555
$ws.Range("I113").Value = 1050.25
$ws.Range("J113").Value = 1934.6
$ws.Range("K113").Value = 1050.25
$ws.Range("L113").Value = 1934.6
$ws.Range("M113").Value = 1119.75
$ws.Range("N113").Value = -6274.6

$ws.Range("H122").Value = 56278588
$ws.Range("I122").Value = 67533840
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 202601520
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -202599070
$ws.Range("N122").Value = -11899.9999

$ws.Range("H133").Value = 39721.05
$ws.Range("J133").Value = 39721.05
$ws.Range("L133").Value = 39721.05
$ws.Range("N133").Value = -49841.05

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 775.125
$ws.Range("I61").Value = 800.1667
$ws.Range("J61").Value = 700
$ws.Range("K61").Value = 800.1667
$ws.Range("L61").Value = 700
$ws.Range("M61").Value = -598.1667
$ws.Range("N61").Value = -1104

$ws.Range("H113").Value = 775.125
$ws.Range("I113").Value = 800.1667
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 800.1667
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 1369.8333
$ws.Range("N113").Value = -5040

$ws.Range("H122").Value = 7712
$ws.Range("I122").Value = 8301
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 24903
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -22453
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3357.3572
$ws.Range("I96").Value = 3667.6667
$ws.Range("J96").Value = 3272.7273
$ws.Range("K96").Value = 3667.6667
$ws.Range("L96").Value = 3272.7273
$ws.Range("M96").Value = -2294.6667
$ws.Range("N96").Value = -6018.7273

$ws.Range("H113").Value = 553.625
$ws.Range("I113").Value = 585.8
$ws.Range("K113").Value = 1757.4
$ws.Range("M113").Value = 412.6000000000001

$ws.Range("H122").Value = 1058.579
$ws.Range("I122").Value = 1058.579
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3175.737
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -725.7370000000001
$ws.Range("N122").ClearContents()
